$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 100
$ws.Range("I13").Value = 100
$ws.Range("K13").Value = 100
$ws.Range("M13").Value = 69
$ws.Range("H29").Value = 549.75
$ws.Range("I29").Value = 99.5
$ws.Range("K29").Value = 298.5
$ws.Range("M29").Value = -17.5
$ws.Range("H43").Value = 1868.25
$ws.Range("I43").Value = 1199.5
$ws.Range("J43").Value = 2002
$ws.Range("K43").Value = 1199.5
$ws.Range("L43").Value = 2002
$ws.Range("M43").Value = -1130.5
$ws.Range("N43").Value = -2140
$ws.Range("H45").Value = 5000
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 5000
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 15000
$ws.Range("N45").Value = -15384
$ws.Range("H52").Value = 808.8
$ws.Range("I52").Value = 227
$ws.Range("J52").Value = 1196.6666
$ws.Range("K52").Value = 681
$ws.Range("L52").Value = 3589.9998
$ws.Range("M52").Value = -521
$ws.Range("N52").Value = -3909.9998
$ws.Range("H59").Value = 996.125
$ws.Range("J59").Value = 1095.5714
$ws.Range("L59").Value = 3286.7142
$ws.Range("N59").Value = -4400.7142
$ws.Range("H69").Value = 3466.6667
$ws.Range("I69").Value = 2156.5
$ws.Range("J69").Value = 4121.75
$ws.Range("K69").Value = 6469.5
$ws.Range("L69").Value = 12365.25
$ws.Range("M69").Value = -5595.5
$ws.Range("N69").Value = -14113.25
$ws.Range("H72").Value = 3466.6667
$ws.Range("I72").Value = 2156.5
$ws.Range("J72").Value = 4121.75
$ws.Range("K72").Value = 19408.5
$ws.Range("L72").Value = 37095.75
$ws.Range("M72").Value = -15040.5
$ws.Range("N72").Value = -45831.75
$ws.Range("H98").Value = 5853864.5
$ws.Range("I98").Value = 6239.3887
$ws.Range("J98").Value = 111111110
$ws.Range("K98").Value = 6239.3887
$ws.Range("L98").Value = 111111110
$ws.Range("M98").Value = -4741.3887
$ws.Range("N98").Value = -111114106
$ws.Range("H122").Value = 5853864.5
$ws.Range("I122").Value = 6239.3887
$ws.Range("J122").Value = 111111110
$ws.Range("K122").Value = 18718.1661
$ws.Range("L122").Value = 333333330
$ws.Range("M122").Value = -16268.1661
$ws.Range("N122").Value = -333338230
$ws.Range("H132").Value = 1254.7727
$ws.Range("I132").Value = 1208.6
$ws.Range("J132").Value = 1716.5
$ws.Range("K132").Value = 3625.8
$ws.Range("L132").Value = 5149.5
$ws.Range("M132").Value = -1095.8
$ws.Range("N132").Value = -10209.5
$ws.Range("H135").Value = 3011.889
$ws.Range("J135").Value = 11250
$ws.Range("L135").Value = 101250
$ws.Range("N135").Value = -106320
$ws.Range("H137").Value = 1563.4
$ws.Range("I137").Value = 1145.6
$ws.Range("J137").Value = 1841.9333
$ws.Range("K137").Value = 3436.8
$ws.Range("L137").Value = 5525.7999
$ws.Range("M137").Value = -886.7999999999997
$ws.Range("N137").Value = -10625.7999
$ws.Range("H139").Value = 78075
$ws.Range("J139").Value = 78075
$ws.Range("L139").Value = 78075
$ws.Range("N139").Value = -88355
$ws.Range("H140").Value = 97300
$ws.Range("J140").Value = 97300
$ws.Range("L140").Value = 97300
$ws.Range("N140").Value = -107660
$ws.Range("M45").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10887065
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 10887065
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 10887065
$ws.Range("N32").Value = -10887639
$ws.Range("H61").Value = 2418.5186
$ws.Range("I61").Value = 2466.125
$ws.Range("J61").Value = 2037.6666
$ws.Range("K61").Value = 2466.125
$ws.Range("L61").Value = 2037.6666
$ws.Range("M61").Value = -2254.125
$ws.Range("N61").Value = -2461.6666
$ws.Range("H74").Value = 1130.125
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("H77").Value = 1130.125
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("H136").Value = 2418.5186
$ws.Range("I136").Value = 2466.125
$ws.Range("J136").Value = 2037.6666
$ws.Range("K136").Value = 7398.375
$ws.Range("L136").Value = 6112.9998
$ws.Range("M136").Value = -4848.375
$ws.Range("N136").Value = -11212.9998
$ws.Range("H141").Value = 65729.89999999999
$ws.Range("J141").Value = 65729.89999999999
$ws.Range("L141").Value = 65729.89999999999
$ws.Range("N141").Value = -76089.89999999999
$ws.Range("M32").ClearContents()
$ws.Range("N74").ClearContents()
$ws.Range("N77").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1126
$ws.Range("I134").Value = 1034.125
$ws.Range("J134").Value = 1371
$ws.Range("K134").Value = 3102.375
$ws.Range("L134").Value = 4113
$ws.Range("M134").Value = -567.375
$ws.Range("N134").Value = -9183

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 15422.232
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 15422.232
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 15422.232
$ws.Range("N31").Value = -16012.232
$ws.Range("H34").Value = 15422.232
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 15422.232
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 15422.232
$ws.Range("N34").Value = -15826.232
$ws.Range("H68").Value = 15835.167
$ws.Range("J68").Value = 15835.167
$ws.Range("L68").Value = 15835.167
$ws.Range("N68").Value = -17333.167
$ws.Range("H71").Value = 15835.167
$ws.Range("J71").Value = 15835.167
$ws.Range("L71").Value = 47505.501
$ws.Range("N71").Value = -54993.501
$ws.Range("H74").Value = 17839.5
$ws.Range("J74").Value = 17839.5
$ws.Range("L74").Value = 17839.5
$ws.Range("N74").Value = -19587.5
$ws.Range("H77").Value = 17839.5
$ws.Range("J77").Value = 17839.5
$ws.Range("L77").Value = 53518.5
$ws.Range("N77").Value = -62254.5
$ws.Range("H99").Value = 1891.4286
$ws.Range("I99").Value = 1936
$ws.Range("J99").Value = 1000
$ws.Range("K99").Value = 1936
$ws.Range("L99").Value = 1000
$ws.Range("M99").Value = -438
$ws.Range("N99").Value = -3996
$ws.Range("H126").Value = 1891.4286
$ws.Range("I126").Value = 1936
$ws.Range("J126").Value = 1000
$ws.Range("K126").Value = 5808
$ws.Range("L126").Value = 3000
$ws.Range("M126").Value = -3338
$ws.Range("N126").Value = -7940
$ws.Range("M31").ClearContents()
$ws.Range("M34").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 997.11536
$ws.Range("I5").Value = 849.0476
$ws.Range("J5").Value = 1619
$ws.Range("K5").Value = 2547.1428
$ws.Range("L5").Value = 4857
$ws.Range("M5").Value = -2435.1428
$ws.Range("N5").Value = -5081
$ws.Range("H38").Value = 361.38095
$ws.Range("I38").Value = 862.5
$ws.Range("J38").Value = 160.93333
$ws.Range("K38").Value = 2587.5
$ws.Range("L38").Value = 482.79999
$ws.Range("M38").Value = -2240.5
$ws.Range("N38").Value = -1176.79999
$ws.Range("H39").Value = 3200
$ws.Range("J39").Value = 6000
$ws.Range("L39").Value = 18000
$ws.Range("N39").Value = -18588
$ws.Range("H131").Value = 9260363
$ws.Range("I131").Value = 1509.091
$ws.Range("J131").Value = 11628906
$ws.Range("K131").Value = 4527.272999999999
$ws.Range("L131").Value = 34886718
$ws.Range("M131").Value = 512.7270000000008
$ws.Range("N131").Value = -34896798
$ws.Range("H135").Value = 997.11536
$ws.Range("I135").Value = 849.0476
$ws.Range("J135").Value = 1619
$ws.Range("K135").Value = 7641.4284
$ws.Range("L135").Value = 14571
$ws.Range("M135").Value = -5106.4284
$ws.Range("N135").Value = -19641

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H138").Value = 69480
$ws.Range("J138").Value = 69480
$ws.Range("L138").Value = 69480
$ws.Range("N138").Value = -79760

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 843.3333
$ws.Range("I22").Value = 533.3333
$ws.Range("J22").Value = 998.3333
$ws.Range("K22").Value = 533.3333
$ws.Range("L22").Value = 998.3333
$ws.Range("M22").Value = -238.3333
$ws.Range("N22").Value = -1588.3333
$ws.Range("H27").Value = 843.3333
$ws.Range("I27").Value = 533.3333
$ws.Range("J27").Value = 998.3333
$ws.Range("K27").Value = 533.3333
$ws.Range("L27").Value = 998.3333
$ws.Range("M27").Value = -426.3333
$ws.Range("N27").Value = -1212.3333
$ws.Range("H40").Value = 11114159
$ws.Range("I40").Value = 3153.6
$ws.Range("J40").Value = 66669188
$ws.Range("K40").Value = 3153.6
$ws.Range("L40").Value = 66669188
$ws.Range("M40").Value = -3017.6
$ws.Range("N40").Value = -66669460
$ws.Range("H46").Value = 8565.308000000001
$ws.Range("I46").Value = 721.2857
$ws.Range("J46").Value = 17716.666
$ws.Range("K46").Value = 721.2857
$ws.Range("L46").Value = 17716.666
$ws.Range("M46").Value = -533.2857
$ws.Range("N46").Value = -18092.666
$ws.Range("H68").Value = 2012
$ws.Range("I68").Value = 1950
$ws.Range("J68").Value = 2148.4
$ws.Range("K68").Value = 1950
$ws.Range("L68").Value = 2148.4
$ws.Range("M68").Value = -1201
$ws.Range("N68").Value = -3646.4
$ws.Range("H71").Value = 2012
$ws.Range("I71").Value = 1950
$ws.Range("J71").Value = 2148.4
$ws.Range("K71").Value = 9750
$ws.Range("L71").Value = 10742
$ws.Range("M71").Value = -6006
$ws.Range("N71").Value = -18230
$ws.Range("H138").Value = 59177.91
$ws.Range("J138").Value = 59177.91
$ws.Range("L138").Value = 59177.91
$ws.Range("N138").Value = -69457.91
$ws.Range("H139").Value = 60850
$ws.Range("I139").Value = 5000
$ws.Range("J139").Value = 79466.664
$ws.Range("K139").Value = 5000
$ws.Range("L139").Value = 79466.664
$ws.Range("M139").Value = 140
$ws.Range("N139").Value = -89746.664

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 5563666
$ws.Range("J2").Value = 9002.143
$ws.Range("L2").Value = 9002.143
$ws.Range("N2").Value = -9226.143
$ws.Range("H127").Value = 69367.375
$ws.Range("J127").Value = 69367.375
$ws.Range("L127").Value = 69367.375
$ws.Range("N127").Value = -79287.375
$ws.Range("H139").Value = 53930.715
$ws.Range("J139").Value = 53930.715
$ws.Range("L139").Value = 53930.715
$ws.Range("N139").Value = -64210.715
